# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 11.945164432584)
    3 = @(0.6753301551942219, 1.667794583268128, 26.21740644021617, 645.3272768299601, 673.8878080086386)
    4 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    5 = @(0.04763786555579896, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 12.65070551845588)
    6 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 616238.5361209477, 616269.6523076545)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
